$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Date" column before the existing "logintime" column (C),
# shifting logintime/logouttime/workingtime one column to the right (D/E/F).
$ws.Range("C1").EntireColumn.Insert()

$ws.Range("C1").Value = "Date"
$ws.Range("C2").Value = '"28/10/21"'

# Approximate the column widths seen in the authored workbook.
$ws.Columns.Item(2).ColumnWidth = 13.25
$ws.Columns.Item(3).ColumnWidth = 7.76

# The former "workingtime" value cell (now F2) becomes a literal text
# timestamp instead of a computed time-serial number, formatted as
# h:mm:ss AM/PM (numFmtId 19).
$ws.Range("F2").NumberFormat = "h:mm:ss AM/PM"
$ws.Range("F2").Value = "03:30:00AM"

$ws.Range("F6").Select()
